$d = $word.ActiveDocument

$replacements = @(
    @("224×9=2016", "227×2=454"),
    @("954×5=4770", "644×5=3220"),
    @("560×5=2800", "371×5=1855"),
    @("681×9=6129", "970×9=8730"),
    @("217×9=1953", "475×7=3325"),
    @("510×8=4080", "365×4=1460"),
    @("430×8=3440", "740×6=4440"),
    @("599×6=3594", "141×7=987"),
    @("850×9=7650", "870×9=7830"),
    @("102×9=918", "965×9=8685"),
    @("467×9=4203", "882×8=7056"),
    @("432×7=3024", "414×2=828"),
    @("906×8=7248", "785×9=7065"),
    @("889×8=7112", "492×9=4428"),
    @("562×8=4496", "383×3=1149"),
    @("561×9=5049", "246×4=984"),
    @("928×5=4640", "186×9=1674"),
    @("220×3=660", "424×9=3816"),
    @("590×6=3540", "405×8=3240"),
    @("385×6=2310", "852×5=4260"),
    @("730×6=4380", "229×8=1832"),
    @("368×8=2944", "882×6=5292"),
    @("585×9=5265", "809×8=6472"),
    @("389×5=1945", "725×6=4350"),
    @("114×6=684", "356×2=712")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
